# Remove the "culture_collection" attribute/column (AK) from the MIGS water
# template. The column itself (header cell + underlying shared string) is
# removed via a normal column delete, which correctly shifts the remaining
# header cells/shared-string table one position to the left. Excel's
# cell-comment ("note") objects, however, do not get re-threaded by the
# simulated column delete here, so we shift the comment *text* ourselves:
# each column's comment takes on the text that used to belong to the next
# column to the right, and the now-superfluous trailing comment is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 37   # AK - the "culture_collection" column being removed
$endCol   = 114  # DJ - last column that currently carries a comment

# 1. Snapshot the current comment text for every column from the one being
#    deleted through the last commented column (read everything up front so
#    later writes in the same range can't clobber values we still need).
$commentText = @{}
for ($c = $startCol; $c -le $endCol; $c++) {
    $cell = $ws.Cells.Item(15, $c)
    if ($cell.Comment -ne $null) {
        $commentText[$c] = $cell.Comment.Text()
    }
}

# 2. Delete the "culture_collection" column itself. This removes the AK15
#    header cell (and its shared string) and shifts every later column's
#    cell content one position to the left.
$ws.Columns.Item($startCol).Delete()

# 3. Re-home the comment text: column c (now holding what used to be column
#    c+1's data) should carry what used to be column c+1's comment.
for ($c = $startCol; $c -le ($endCol - 1); $c++) {
    $null = $ws.Cells.Item(15, $c).Comment.Text($commentText[$c + 1])
}

# 4. The comment that used to sit on the last column no longer corresponds
#    to any column (everything shifted left by one), so drop it.
$null = $ws.Cells.Item(15, $endCol).Comment.Delete()
